$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.404.02'
$ws.Range('E2').Value = '  +0.58%  '

$ws.Range('D3').Value = '3.843.26'
$ws.Range('E3').Value = '  +1.15%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').Value = "'714.34"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.73%  '

$ws.Range('D6').Value = "'173.08"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').Value = '3.841.47'

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('E9').Value = '  -0.05%  '

$ws.Range('D10').Value = "'0.164"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.52%  '

$ws.Range('D11').Value = "'7.36"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.87%  '

$ws.Range('D12').Value = "'0.462"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.24%  '

$ws.Range('D13').Value = "'0.0000256"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.30%  '

$ws.Range('D14').Value = "'36.88"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.21%  '

$ws.Range('D15').Value = '4.491.48'
$ws.Range('E15').Value = '  +1.14%  '

$ws.Range('D16').Value = '3.864.15'
$ws.Range('E16').Value = '  +1.37%  '

$ws.Range('D17').Value = '71.336.63'
$ws.Range('E17').Value = '  +0.55%  '

$ws.Range('D18').Value = "'7.28"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.04%  '

$ws.Range('E19').Value = '  +0.40%  '

$ws.Range('D20').Value = "'17.47"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.01%  '

$ws.Range('D21').Value = "'500.32"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +4.16%  '

$ws.Range('D22').Value = "'10.72"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.26%  '

$ws.Range('D23').Value = "'0.736"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.30%  '

$ws.Range('D24').Value = "'85.49"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.60%  '

$ws.Range('E25').Value = '  +2.32%  '

$ws.Range('E26').Value = '  +1.54%  '

$ws.Range('D27').Value = "'12.23"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.55%  '

$ws.Range('E28').Value = '  -2.69%  '

$ws.Range('D29').Value = "'3.17"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.92%  '

$ws.Range('E30').Value = '  -0.01%  '

$ws.Range('D31').Value = "'7.51"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.76%  '

$ws.Range('D32').Value = "'2.25"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.90%  '

$ws.Range('D33').Value = "'29.47"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '

$ws.Range('D34').Value = "'0.182"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.20%  '

$ws.Range('D35').Value = "'9.24"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.31%  '

$ws.Range('D36').Value = '3.808.27'
$ws.Range('E36').Value = '  +1.53%  '

$ws.Range('E37').Value = '  -0.14%  '

$ws.Range('E38').Value = '  +0.32%  '

$ws.Range('E39').Value = '  +5.96%  '

$ws.Range('E40').Value = '  +0.47%  '

$ws.Range('D41').Value = "'3.36"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.43%  '

$ws.Range('D42').Value = "'2.29"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.54%  '

$ws.Range('E44').Value = '  +0.07%  '

$ws.Range('E45').Value = '  +0.48%  '

$ws.Range('D46').Value = "'163.58"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.33%  '

$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = "'427.68"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.29%  '

$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = "'49.01"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.27%  '

$ws.Range('D49').Value = "'8.75"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.38%  '

$ws.Range('E50').Value = '  +0.32%  '

$ws.Range('E51').Value = '  -0.87%  '

